$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 2 and row 3 for columns D, M, N, O, P, Q, S, T

# Row 2 -> new values (previously row 3's values)
$ws.Range("D2").Value = 44874
$ws.Range("M2").Value = 67
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1600
$ws.Range("T2").Value = 10

# Row 3 -> new values (previously row 2's values)
$ws.Range("D3").Value = 44855
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("S3").Value = 3000
$ws.Range("T3").Value = 5
